$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B5").Value = [double]"-8.11407673534362e-09"
$ws.Range("C5").Value = [double]"-8.114076735446819e-09"
$ws.Range("D5").Value = [double]"-8.114076735509334e-09"
$ws.Range("E5").Value = [double]"-8.114076735543718e-09"
$ws.Range("F5").Value = [double]"-8.114076735149643e-09"
$ws.Range("B6").Value = 0.04267278849067333
$ws.Range("C6").Value = 0.04289296817092081
$ws.Range("D6").Value = 0.043026375592555566
$ws.Range("E6").Value = 0.04309975312991184
$ws.Range("F6").Value = 0.04225941696976104
$ws.Range("B7").Value = -0.18616990724214555
$ws.Range("C7").Value = -0.18634528991700805
$ws.Range("D7").Value = -0.18645157012503108
$ws.Range("E7").Value = -0.18651002881727677
$ws.Range("F7").Value = -0.18584064589797633
$ws.Range("B8").Value = 0.14349711063739545
$ws.Range("C8").Value = 0.1434523136320105
$ws.Range("D8").Value = 0.14342518641839877
$ws.Range("E8").Value = 0.14341026757328817
$ws.Range("F8").Value = 0.14358122081413854
$ws.Range("B9").Value = 0.005000571413263254
$ws.Range("C9").Value = 0.0051806039760951
$ws.Range("D9").Value = 0.0052472829146299966
$ws.Range("E9").Value = 0.005307294056865159
$ws.Range("F9").Value = 0.005213943428881287
$ws.Range("B10").Value = 0.009998062702596202
$ws.Range("C10").Value = 0.010358011163705293
$ws.Range("D10").Value = 0.010491325447960506
$ws.Range("E10").Value = 0.010611308244979203
$ws.Range("F10").Value = 0.010424667425185137
$ws.Range("B11").Value = 0.017669841922028157
$ws.Range("C11").Value = 0.01699000255435733
$ws.Range("D11").Value = 0.01679006037467482
$ws.Range("E11").Value = 0.016563438988404514
$ws.Range("F11").Value = 0.016189974347630958
$ws.Range("B12").Value = 0.005000668208116438
$ws.Range("C12").Value = 0.0051807078634956375
$ws.Range("D12").Value = 0.00524738949265046
$ws.Range("E12").Value = 0.005307403085154682
$ws.Range("F12").Value = 0.005214048649905797
$ws.Range("B13").Value = 0.005000630883280346
$ws.Range("C13").Value = 0.005180667802565203
$ws.Range("D13").Value = 0.005247348393776196
$ws.Range("E13").Value = 0.005307361040859895
$ws.Range("F13").Value = 0.005214008072517727
$ws.Range("B14").Value = 0.005000309305316134
$ws.Range("C14").Value = 0.005180322652479951
$ws.Range("D14").Value = 0.0052469943015812604
$ws.Range("E14").Value = 0.005306998803677514
$ws.Range("F14").Value = 0.005213658473036318
$ws.Range("B15").Value = 0.02398164072690172
$ws.Range("C15").Value = 0.024521582990941784
$ws.Range("D15").Value = 0.024691603687843535
$ws.Range("E15").Value = 0.024812503388652454
$ws.Range("F15").Value = 0.02441143085231938
$ws.Range("B16").Value = -1.4221559355755569
$ws.Range("C16").Value = -1.4221137819375027
$ws.Range("D16").Value = -1.4220882280638245
$ws.Range("E16").Value = -1.4220741693655379
$ws.Range("F16").Value = -1.4222350234058094
$ws.Range("B17").Value = -1.512188508422315
$ws.Range("C17").Value = -1.5121838355430783
$ws.Range("D17").Value = -1.5121810183281013
$ws.Range("E17").Value = -1.5121794625517695
$ws.Range("F17").Value = -1.512197167428314
$ws.Range("B18").Value = 0.00250007601218117
$ws.Range("C18").Value = 0.0025900909121600534
$ws.Range("D18").Value = 0.002623429978591316
$ws.Range("E18").Value = 0.0026534352373290483
$ws.Range("F18").Value = 0.002606760428658602
$ws.Range("B19").Value = 0.005486378506089839
$ws.Range("C19").Value = 0.0056526226661973615
$ws.Range("D19").Value = 0.005714217264325665
$ws.Range("E19").Value = 0.005769853927509638
$ws.Range("F19").Value = 0.005685351023857658
$ws.Range("B20").Value = 0.010980998086839227
$ws.Range("C20").Value = 0.011084410356678278
$ws.Range("D20").Value = 0.011092712039849965
$ws.Range("E20").Value = 0.011067895241514823
$ws.Range("F20").Value = 0.0108918705180177
$ws.Range("B21").Value = 0.0025001738319419646
$ws.Range("C21").Value = 0.0025901958995573145
$ws.Range("D21").Value = 0.002623537685094416
$ws.Range("E21").Value = 0.0026535454200406258
$ws.Range("F21").Value = 0.0026068667637798356
$ws.Range("B22").Value = 0.00250019785738736
$ws.Range("C22").Value = 0.0025902216859692846
$ws.Range("D22").Value = 0.0026235641395274466
$ws.Range("E22").Value = 0.0026535724829439837
$ws.Range("F22").Value = 0.002606892882584602
$ws.Range("B23").Value = 0.0024998716781486205
$ws.Range("C23").Value = 0.0025898715973295544
$ws.Range("D23").Value = 0.0026232049808305163
$ws.Range("E23").Value = 0.0026532050627176575
$ws.Range("F23").Value = 0.0026065382808848577
$ws.Range("B24").Value = 0.0022585407879416454
$ws.Range("C24").Value = 0.002300984733434881
$ws.Range("D24").Value = 0.002319675975362818
$ws.Range("E24").Value = 0.0023312743871574657
$ws.Range("F24").Value = 0.0022450525323725206
$ws.Range("B25").Value = 1.9390958049051548
$ws.Range("C25").Value = 1.9390443734749772
$ws.Range("D25").Value = 1.9390132117717318
$ws.Range("E25").Value = 1.93899607314067
$ws.Range("F25").Value = 1.9391923908364115
$ws.Range("B26").Value = 2.174856730854773
$ws.Range("C26").Value = 2.1748312930174447
$ws.Range("D26").Value = 2.1748159263729168
$ws.Range("E26").Value = 2.174807457636802
$ws.Range("F26").Value = 2.1749041758382623
$ws.Range("B27").Value = [double]"2.9477240831714792e-05"
$ws.Range("C27").Value = [double]"3.163821969727627e-05"
$ws.Range("D27").Value = [double]"3.24579847037672e-05"
$ws.Range("E27").Value = [double]"3.320471999513211e-05"
$ws.Range("F27").Value = [double]"3.204659068193962e-05"
$ws.Range("B28").Value = 0.00012265950077622456
$ws.Range("C28").Value = 0.0001313170285372922
$ws.Range("D28").Value = 0.00013459903746485555
$ws.Range("E28").Value = 0.0001375897453592473
$ws.Range("F28").Value = 0.00013297280163801077
$ws.Range("B29").Value = 0.00040817369267680693
$ws.Range("C29").Value = 0.000388106490190614
$ws.Range("D29").Value = 0.00038191144746627467
$ws.Range("E29").Value = 0.0003742649568761022
$ws.Range("F29").Value = 0.00035907983395989663
$ws.Range("B30").Value = [double]"2.949258464304242e-05"
$ws.Range("C30").Value = [double]"3.165444217345958e-05"
$ws.Range("D30").Value = [double]"3.247454541748751e-05"
$ws.Range("E30").Value = [double]"3.3221542701979445e-05"
$ws.Range("F30").Value = [double]"3.2062507403936046e-05"
$ws.Range("B31").Value = [double]"2.9492345777144864e-05"
$ws.Range("C31").Value = [double]"3.165417656709492e-05"
$ws.Range("D31").Value = [double]"3.247426942214567e-05"
$ws.Range("E31").Value = [double]"3.3221257129081806e-05"
$ws.Range("F31").Value = [double]"3.206223664770651e-05"
$ws.Range("B32").Value = [double]"2.9492345777144857e-05"
$ws.Range("C32").Value = [double]"3.1654176567094915e-05"
$ws.Range("D32").Value = [double]"3.247426942214569e-05"
$ws.Range("E32").Value = [double]"3.32212571290818e-05"
$ws.Range("F32").Value = [double]"3.2062236647706516e-05"
$ws.Range("B34").Value = 1.0603466825226278
$ws.Range("C34").Value = 1.0603386919800197
$ws.Range("D34").Value = 1.0603357895706005
$ws.Range("E34").Value = 1.060333855243991
$ws.Range("F34").Value = 1.0603438708737556
$ws.Range("B35").Value = 1.0859808453709558
$ws.Range("C35").Value = 1.0859815575394023
$ws.Range("D35").Value = 1.0859819889559348
$ws.Range("E35").Value = 1.0859822262048857
$ws.Range("F35").Value = 1.0859795075184864
$ws.Range("B36").Value = 1.0601823387545997
$ws.Range("C36").Value = 1.0601684315710587
$ws.Range("D36").Value = 1.0601633378084658
$ws.Range("E36").Value = 1.0601594312618416
$ws.Range("F36").Value = 1.0601725147885723
$ws.Range("B37").Value = 1.060218657731482
$ws.Range("C37").Value = 1.0602064442835548
$ws.Range("D37").Value = 1.0602019775584375
$ws.Range("E37").Value = 1.060198632862457
$ws.Range("F37").Value = 1.0602108171911862
$ws.Range("B38").Value = 1.0598444353356313
$ws.Range("C38").Value = 1.0598464689767166
$ws.Range("D38").Value = 1.059846905636889
$ws.Range("E38").Value = 1.0598493610120587
$ws.Range("F38").Value = 1.0598689996251651
$ws.Range("B39").Value = 1.0596307683183084
$ws.Range("C39").Value = 1.0596251094983642
$ws.Range("D39").Value = 1.0596226970909828
$ws.Range("E39").Value = 1.0596225883017243
$ws.Range("F39").Value = 1.0596462156138453
$ws.Range("B40").Value = 1.0596800875791859
$ws.Range("C40").Value = 1.0596762042870431
$ws.Range("D40").Value = 1.059674449483158
$ws.Range("E40").Value = 1.0596749325373254
$ws.Range("F40").Value = 1.0596976392042103
$ws.Range("B41").Value = 1.059639350533384
$ws.Range("C41").Value = 1.0596340006844183
$ws.Range("D41").Value = 1.0596317027109308
$ws.Range("E41").Value = 1.0596316969122677
$ws.Range("F41").Value = 1.059655164016826
$ws.Range("B42").Value = 1.0823989433572518
$ws.Range("C42").Value = 1.082399552851273
$ws.Range("D42").Value = 1.082399922133606
$ws.Range("E42").Value = 1.0824001251219446
$ws.Range("F42").Value = 1.0823977971651118
$ws.Range("B43").Value = -0.0008919103824805479
$ws.Range("C43").Value = -0.0006316482174238614
$ws.Range("D43").Value = -0.0005055935395941902
$ws.Range("E43").Value = -0.00035964865734512333
$ws.Range("F43").Value = -0.00036823806952544563
$ws.Range("B44").Value = [double]"9.967310392012425e-22"
$ws.Range("C44").Value = [double]"-9.829717723932946e-20"
$ws.Range("D44").Value = [double]"5.999115317587511e-20"
$ws.Range("E44").Value = [double]"4.093159761592188e-20"
$ws.Range("F44").Value = [double]"4.093159761592188e-20"
$ws.Range("B45").Value = -0.00048537153810864213
$ws.Range("C45").Value = -0.00047155327005715377
$ws.Range("D45").Value = -0.0004664576139397097
$ws.Range("E45").Value = -0.00046207282058258284
$ws.Range("F45").Value = -0.00047093640945295893
$ws.Range("B46").Value = [double]"-8.485739208970507e-05"
$ws.Range("C46").Value = [double]"-8.86543429255147e-05"
$ws.Range("D46").Value = [double]"-8.971123458695566e-05"
$ws.Range("E46").Value = [double]"-9.088784199941676e-05"
$ws.Range("F46").Value = [double]"-9.310421310054262e-05"
$ws.Range("B47").Value = [double]"2.0993320770248966e-06"
$ws.Range("C47").Value = [double]"2.106139472085875e-06"
$ws.Range("D47").Value = [double]"2.108555783440064e-06"
$ws.Range("E47").Value = [double]"2.1106361979856894e-06"
$ws.Range("F47").Value = [double]"2.1068380473971538e-06"
$ws.Range("B48").Value = [double]"2.374521956174667e-06"
$ws.Range("C48").Value = [double]"2.3739899784138335e-06"
$ws.Range("D48").Value = [double]"2.3737793004399404e-06"
$ws.Range("E48").Value = [double]"2.3732545600553326e-06"
$ws.Range("F48").Value = [double]"2.370592455300252e-06"
$ws.Range("B49").Value = [double]"1.906608213206597e-06"
$ws.Range("C49").Value = [double]"1.9010682622889227e-06"
$ws.Range("D49").Value = [double]"1.8990998655082748e-06"
$ws.Range("E49").Value = [double]"1.8974041550411791e-06"
$ws.Range("F49").Value = [double]"1.900499678040642e-06"
$ws.Range("B50").Value = [double]"2.140784145672382e-06"
$ws.Range("C50").Value = [double]"2.1412174339978615e-06"
$ws.Range("D50").Value = [double]"2.1413887907442063e-06"
$ws.Range("E50").Value = [double]"2.141815915988925e-06"
$ws.Range("F50").Value = [double]"2.1439834845688995e-06"
$ws.Range("B51").Value = 0.011663737409459272
$ws.Range("C51").Value = 0.011663737117471134
$ws.Range("D51").Value = 0.011663741192975394
$ws.Range("E51").Value = 0.011663749029831859
$ws.Range("F51").Value = 0.01166375
$ws.Range("B52").Value = 0.0058336273388629204
$ws.Range("C52").Value = 0.005833628716253995
$ws.Range("D52").Value = 0.005833629713125472
$ws.Range("E52").Value = 0.005833629761888063
$ws.Range("F52").Value = 0.005833624999999999
